# Apply cryptos-list price/volume refresh per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as literal text
# (mirrors typing a leading apostrophe in Excel), so e.g. "1.00" keeps
# its trailing zero instead of Excel coercing it to the number 1.
function Set-TextValue($rangeAddress, $text) {
    $ws.Range($rangeAddress).Value = "'" + $text
}

$ws.Range("D2").Value = '79.514.32'
$ws.Range("E2").Value = '  +4.24%  '

$ws.Range("D3").Value = '3.156.97'
$ws.Range("E3").Value = '  +2.61%  '

Set-TextValue "D4" '1.00'
$ws.Range("E4").Value = '  +0.27%  '

Set-TextValue "D5" '207.44'
$ws.Range("E5").Value = '  +4.94%  '

Set-TextValue "D6" '627.16'
$ws.Range("E6").Value = '  +1.97%  '

Set-TextValue "D7" '0.269'
$ws.Range("E7").Value = '  +28.90%  '

Set-TextValue "D8" '1.00'
$ws.Range("E8").Value = '  +0.09%  '

Set-TextValue "D9" '0.593'
$ws.Range("E9").Value = '  +7.84%  '

$ws.Range("D10").Value = '3.159.33'
$ws.Range("E10").Value = '  +2.88%  '

Set-TextValue "D11" '0.603'
$ws.Range("E11").Value = '  +36.85%  '

Set-TextValue "D12" '0.0000253'
$ws.Range("E12").Value = '  +31.34%  '

$ws.Range("E13").Value = '  +2.42%  '

Set-TextValue "D14" '5.28'
$ws.Range("E14").Value = '  +1.50%  '

$ws.Range("D15").Value = '3.742.34'
$ws.Range("E15").Value = '  +2.76%  '

Set-TextValue "D16" '31.58'
$ws.Range("E16").Value = '  +8.59%  '

$ws.Range("D17").Value = '79.633.86'
$ws.Range("E17").Value = '  +4.41%  '

$ws.Range("D18").Value = '3.163.75'
$ws.Range("E18").Value = '  +2.58%  '

Set-TextValue "D19" '14.34'
$ws.Range("E19").Value = '  +5.69%  '

Set-TextValue "D20" '438.94'
$ws.Range("E20").Value = '  +15.66%  '

$ws.Range("E21").Value = '  +19.00%  '

Set-TextValue "D22" '9.17'
$ws.Range("E22").Value = '  +0.36%  '

Set-TextValue "D23" '5.25'
$ws.Range("E23").Value = '  +19.28%  '

Set-TextValue "D24" '6.79'
$ws.Range("E24").Value = '  +4.94%  '

$ws.Range("D25").Value = '3.330.80'
$ws.Range("E25").Value = '  +2.92%  '

Set-TextValue "D26" '76.16'
$ws.Range("E26").Value = '  +5.55%  '

Set-TextValue "D27" '4.70'
$ws.Range("E27").Value = '  +8.27%  '

$ws.Range("E28").Value = '  +10.70%  '

Set-TextValue "D29" '0.998'
$ws.Range("E29").Value = '  -0.22%  '

Set-TextValue "D30" '0.0000122'
$ws.Range("E30").Value = '  +13.37%  '

Set-TextValue "D31" '9.07'
$ws.Range("E31").Value = '  +9.29%  '

Set-TextValue "D32" '0.998'
$ws.Range("E32").Value = '  +0.31%  '

Set-TextValue "D33" '549.81'
$ws.Range("E33").Value = '  +10.52%  '

Set-TextValue "D34" '1.47'
$ws.Range("E34").Value = '  +4.47%  '

$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D35" '0.150'
$ws.Range("E35").Value = '  +21.77%  '

$ws.Range("B36").Value = 'PancakeSwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D36" '1.99'
$ws.Range("E36").Value = '  +4.36%  '

Set-TextValue "D37" '23.19'
$ws.Range("E37").Value = '  +12.00%  '

Set-TextValue "D38" '0.120'
$ws.Range("E38").Value = '  +17.52%  '

$ws.Range("E39").Value = '  +0.01%  '

Set-TextValue "D40" '0.407'
$ws.Range("E40").Value = '  +7.76%  '

$ws.Range("E41").Value = '  +3.57%  '

Set-TextValue "D42" '164.25'
$ws.Range("E42").Value = '  +1.56%  '

Set-TextValue "D43" '5.65'
$ws.Range("E43").Value = '  +10.81%  '

Set-TextValue "D45" '188.51'
$ws.Range("E45").Value = '  -2.89%  '

Set-TextValue "D46" '1.81'
$ws.Range("E46").Value = '  +10.56%  '

Set-TextValue "D47" '2.68'
$ws.Range("E47").Value = '  +10.78%  '

Set-TextValue "D48" '0.783'
$ws.Range("E48").Value = '  -1.96%  '

$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D49" '43.50'
$ws.Range("E49").Value = '  +5.28%  '

$ws.Range("B50").Value = 'ImmutableX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D50" '1.31'
$ws.Range("E50").Value = '  +5.06%  '

Set-TextValue "D51" '4.26'
$ws.Range("E51").Value = '  +9.93%  '
